# Generate Report for Handback
# Updates row 6 ("14368ccb-8ec5-4156-9fb1-9c3b7a22bce6.md") on the zh-cn and
# de-de sheets with the handback result: target file, handback xliff file,
# handback datetime, and an error detail because the handback file version
# is stale. Also widens the "Error Detail" column and adds a hyperlink on
# the newly-populated "Latest Target File" cell.

$wb = $excel.ActiveWorkbook

$latestTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bd5c62fec573272d92ab76fb96fa8ccf3e03dd63/e2e/14368ccb-8ec5-4156-9fb1-9c3b7a22bce6.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e99e13cbca8d8c90552666f5a748717d14e4981/e2e/14368ccb-8ec5-4156-9fb1-9c3b7a22bce6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bd5c62fec573272d92ab76fb96fa8ccf3e03dd63/e2e/14368ccb-8ec5-4156-9fb1-9c3b7a22bce6.md."

$configs = @(
    @{
        Sheet = "zh-cn"
        TargetFileName = "14368ccb-8ec5-4156-9fb1-9c3b7a22bce6.md"
        HandbackXliff = "14368ccb-8ec5-4156-9fb1-9c3b7a22bce6.ffa07c3b9fe6e39dcbb578d364c82c2612c9cf97.zh-cn.xlf"
        HandbackDateTime = "2016-08-27 02:42:26"
    },
    @{
        Sheet = "de-de"
        TargetFileName = "14368ccb-8ec5-4156-9fb1-9c3b7a22bce6.md"
        HandbackXliff = "14368ccb-8ec5-4156-9fb1-9c3b7a22bce6.ffa07c3b9fe6e39dcbb578d364c82c2612c9cf97.de-de.xlf"
        HandbackDateTime = "2016-08-27 02:42:33"
    }
)

foreach ($cfg in $configs) {
    $ws = $wb.Worksheets.Item($cfg.Sheet)

    # Widen column P (Error Detail, column 16) to 40 characters.
    $ws.Columns.Item(16).ColumnWidth = 39.166666666666664

    # Populate the newly-arrived handback data on row 6.
    $ws.Range("I6").Value = $cfg.TargetFileName
    $ws.Range("J6").Value = $cfg.HandbackXliff
    $ws.Range("K6").Value = $cfg.HandbackDateTime
    $ws.Range("P6").Value = $errorDetail

    # Style I6 like the other "HyperLink" styled cells (underline + blue).
    $ws.Range("I6").Font.Underline = 2
    $ws.Range("I6").Font.Color = 15570276

    # Rebuild the hyperlink list so the new I6 hyperlink lands in the same
    # position (right after A6) as in the canonical workbook, with every
    # following hyperlink's relationship id shifting up by one.
    $addrs = New-Object System.Collections.ArrayList
    $rows = New-Object System.Collections.ArrayList
    $cols = New-Object System.Collections.ArrayList
    $disps = New-Object System.Collections.ArrayList

    for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
        $h = $ws.Hyperlinks.Item($i)
        [void]$addrs.Add($h.Address)
        [void]$rows.Add($h.Range.Row)
        [void]$cols.Add($h.Range.Column)
        [void]$disps.Add($h.TextToDisplay)
    }

    # Clear every hyperlink on the sheet (this also drops them from the
    # relationship file) so we can re-add them in the desired order.
    $ws.Range("A1").Hyperlinks.Delete()

    for ($i = 0; $i -lt $addrs.Count; $i++) {
        if ($rows[$i] -lt 7) {
            $ws.Hyperlinks.Add($ws.Cells.Item($rows[$i], $cols[$i]), $addrs[$i], "", "", $disps[$i]) | Out-Null
        }
    }

    $ws.Hyperlinks.Add($ws.Range("I6"), $latestTargetUrl, "", "", $cfg.TargetFileName) | Out-Null

    for ($i = 0; $i -lt $addrs.Count; $i++) {
        if ($rows[$i] -ge 7) {
            $ws.Hyperlinks.Add($ws.Cells.Item($rows[$i], $cols[$i]), $addrs[$i], "", "", $disps[$i]) | Out-Null
        }
    }
}

Write-Host "Report generated for handback."
